$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 ("runs"/"balls" for Tom Curran vs Rajasthan Royals) and Row 4 values
# are swapped: C2/D2 <- former C4/D4 values, C4/D4 <- former C2/D2 values.
# Values are entered with a leading apostrophe so they stay text (matching
# the text-as-number storage used by the rest of the sheet) instead of
# being converted to numeric cells.
$ws.Range("C2").Value = "'4"
$ws.Range("D2").Value = "'1"
$ws.Range("C4").Value = "'15"
$ws.Range("D4").Value = "'16"
